$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 309 (existing rows 309:345 shift down to 312:348).
$ws.Range("A309:A311").EntireRow.Insert()

# New record date shared by the 3 inserted rows: 2023-05-31 (serial 45077).
$newDate = Get-Date -Year 2023 -Month 5 -Day 31 -Hour 0 -Minute 0 -Second 0

# Columns A:J are identical for every row in this block.
$aj = New-Object 'object[,]' 3,10
for ($i = 0; $i -lt 3; $i++) {
    $aj[$i,0] = 7
    $aj[$i,1] = "Terminal Hortofrutícola Agro Chillán"
    $aj[$i,2] = "Ñuble"
    $aj[$i,3] = $newDate
    $aj[$i,4] = 16
    $aj[$i,5] = "Fruta"
    $aj[$i,6] = 100104
    $aj[$i,7] = "Frutos de pepita"
    $aj[$i,8] = 100104005
    $aj[$i,9] = "Pera"
}
$ws.Range("A309:J311").Value = $aj

# Columns K:T hold the new "Winter Nelis" records.
$kt = New-Object 'object[,]' 3,10
$kt[0,0] = "Winter Nelis"; $kt[0,1] = "Especial"; $kt[0,2] = 60; $kt[0,3] = 12000; $kt[0,4] = 12000; $kt[0,5] = 12000; $kt[0,6] = "$/bandeja 18 kilos granel"; $kt[0,7] = "Región de O'Higgins"; $kt[0,8] = 667; $kt[0,9] = 18
$kt[1,0] = "Winter Nelis"; $kt[1,1] = "Primera";  $kt[1,2] = 60; $kt[1,3] = 10000; $kt[1,4] = 10000; $kt[1,5] = 10000; $kt[1,6] = "$/bandeja 18 kilos granel"; $kt[1,7] = "Región de O'Higgins"; $kt[1,8] = 556; $kt[1,9] = 18
$kt[2,0] = "Winter Nelis"; $kt[2,1] = "Segunda";  $kt[2,2] = 40; $kt[2,3] = 8000;  $kt[2,4] = 8000;  $kt[2,5] = 8000;  $kt[2,6] = "$/bandeja 18 kilos granel"; $kt[2,7] = "Región de O'Higgins"; $kt[2,8] = 444; $kt[2,9] = 18
$ws.Range("K309:T311").Value = $kt
